$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing data rows (2-6), keep header row 1 intact
$ws.Rows("2:6").Delete()

# Row 2
$ws.Cells.Item(2,1).Value = 12579
$ws.Cells.Item(2,2).Value = "Elijah Henson"
$ws.Cells.Item(2,3).Value = 38
$ws.Cells.Item(2,4).Value = "M"
$ws.Cells.Item(2,5).Value = "Mexe"
$ws.Cells.Item(2,6).Value = 44380
$ws.Cells.Item(2,7).Value = 44381
$ws.Cells.Item(2,8).Value = 44387
$ws.Cells.Item(2,9).Value = 44388
$ws.Cells.Item(2,10).Value = 44379
$ws.Cells.Item(2,11).Value = 44382

# Row 3
$ws.Cells.Item(3,1).Value = 13289
$ws.Cells.Item(3,2).Value = "Ella-Mai Gregory"
$ws.Cells.Item(3,3).Value = 58
$ws.Cells.Item(3,4).Value = "M"
$ws.Cells.Item(3,5).Value = "Grand Wellworth"
$ws.Cells.Item(3,6).Value = 44379
$ws.Cells.Item(3,7).Value = 44380
$ws.Cells.Item(3,10).Value = 44379
$ws.Cells.Item(3,11).Value = 44381

# Row 4
$ws.Cells.Item(4,1).Value = 13479
$ws.Cells.Item(4,2).Value = "Ceara West"
$ws.Cells.Item(4,3).Value = 2
$ws.Cells.Item(4,4).Value = "F"
$ws.Cells.Item(4,5).Value = "Chorgains"
$ws.Cells.Item(4,6).Value = 44378
$ws.Cells.Item(4,7).Value = 44379
$ws.Cells.Item(4,8).Value = 44392
$ws.Cells.Item(4,9).Value = 44393
$ws.Cells.Item(4,10).Value = 44378
$ws.Cells.Item(4,11).Value = 44380

# Row 5
$ws.Cells.Item(5,1).Value = 13547
$ws.Cells.Item(5,2).Value = "Francissek Vickers"
$ws.Cells.Item(5,3).Value = 2
$ws.Cells.Item(5,4).Value = "M"
$ws.Cells.Item(5,5).Value = "Eastmsallbuck Creek"
$ws.Cells.Item(5,6).Value = 44387
$ws.Cells.Item(5,7).Value = 44388
$ws.Cells.Item(5,8).Value = 44392
$ws.Cells.Item(5,9).Value = 44393
$ws.Cells.Item(5,10).Value = 44387
$ws.Cells.Item(5,11).Value = 44389

# Row 6
$ws.Cells.Item(6,1).Value = 13566
$ws.Cells.Item(6,2).Value = "Penelope F. Fields"
$ws.Cells.Item(6,3).Value = 45
$ws.Cells.Item(6,4).Value = "F"
$ws.Cells.Item(6,5).Value = "San Wadhor"
$ws.Cells.Item(6,6).Value = 44379
$ws.Cells.Item(6,7).Value = 44379
$ws.Cells.Item(6,10).Value = 44378
$ws.Cells.Item(6,11).Value = 44380

# Row 7
$ws.Cells.Item(7,1).Value = 13597
$ws.Cells.Item(7,2).Value = "Agata Lucas"
$ws.Cells.Item(7,3).Value = 35
$ws.Cells.Item(7,4).Value = "F"
$ws.Cells.Item(7,5).Value = "Port Sipleach"
$ws.Cells.Item(7,6).Value = 44379
$ws.Cells.Item(7,7).Value = 44380
$ws.Cells.Item(7,10).Value = 44377
$ws.Cells.Item(7,11).Value = 44381

# Row 8
$ws.Cells.Item(8,1).Value = 13788
$ws.Cells.Item(8,2).Value = "Eve M. Mcbride"
$ws.Cells.Item(8,3).Value = 58
$ws.Cells.Item(8,4).Value = "F"
$ws.Cells.Item(8,5).Value = "San Wadhor"
$ws.Cells.Item(8,6).Value = 44379
$ws.Cells.Item(8,7).Value = 44379
$ws.Cells.Item(8,10).Value = 44377
$ws.Cells.Item(8,11).Value = 44380

# Row 9
$ws.Cells.Item(9,1).Value = 18400
$ws.Cells.Item(9,2).Value = "Leonidas Hudson"
$ws.Cells.Item(9,3).Value = 14
$ws.Cells.Item(9,4).Value = "M"
$ws.Cells.Item(9,5).Value = "Eastmsallbuck Creek"
$ws.Cells.Item(9,6).Value = 44384
$ws.Cells.Item(9,7).Value = 44385
$ws.Cells.Item(9,8).Value = 44392
$ws.Cells.Item(9,9).Value = 44392
$ws.Cells.Item(9,10).Value = 44383
$ws.Cells.Item(9,11).Value = 44386

# Row 10
$ws.Cells.Item(10,1).Value = 18793
$ws.Cells.Item(10,2).Value = "Dustin Payne"
$ws.Cells.Item(10,3).Value = 10
$ws.Cells.Item(10,4).Value = "M"
$ws.Cells.Item(10,5).Value = "Grand Wellworth"
$ws.Cells.Item(10,6).Value = 44385
$ws.Cells.Item(10,7).Value = 44386
$ws.Cells.Item(10,10).Value = 44385
$ws.Cells.Item(10,11).Value = 44387

# Row 11
$ws.Cells.Item(11,1).Value = 44980
$ws.Cells.Item(11,2).Value = "Amal Ford"
$ws.Cells.Item(11,3).Value = 40
$ws.Cells.Item(11,4).Value = "M"
$ws.Cells.Item(11,5).Value = "Grand Wellworth"
$ws.Cells.Item(11,6).Value = 44392
$ws.Cells.Item(11,7).Value = 44392
$ws.Cells.Item(11,10).Value = 44387
$ws.Cells.Item(11,11).Value = 44392.66666666667

# Row 12
$ws.Cells.Item(12,1).Value = 44986
$ws.Cells.Item(12,2).Value = "Martin F Romero"
$ws.Cells.Item(12,3).Value = 18
$ws.Cells.Item(12,4).Value = "M"
$ws.Cells.Item(12,5).Value = "Port Sipleach"
$ws.Cells.Item(12,6).Value = 44392
$ws.Cells.Item(12,7).Value = 44392
$ws.Cells.Item(12,10).Value = 44389
$ws.Cells.Item(12,11).Value = 44392.66666666667

# Row 13
$ws.Cells.Item(13,1).Value = 44990
$ws.Cells.Item(13,2).Value = "Fern Christian Mcarthur"
$ws.Cells.Item(13,3).Value = 40
$ws.Cells.Item(13,4).Value = "M"
$ws.Cells.Item(13,5).Value = "Port Sipleach"
$ws.Cells.Item(13,6).Value = 44392
$ws.Cells.Item(13,7).Value = 44392
$ws.Cells.Item(13,10).Value = 44391
$ws.Cells.Item(13,11).Value = 44392.66666666667

# Row 14
$ws.Cells.Item(14,1).Value = 44992
$ws.Cells.Item(14,2).Value = "Jessica Bauer"
$ws.Cells.Item(14,3).Value = 3
$ws.Cells.Item(14,4).Value = "F"
$ws.Cells.Item(14,5).Value = "Eastmsallbuck Creek"
$ws.Cells.Item(14,6).Value = 44392
$ws.Cells.Item(14,7).Value = 44392
$ws.Cells.Item(14,10).Value = 44387
$ws.Cells.Item(14,11).Value = 44392.66666666667

# Row 15
$ws.Cells.Item(15,1).Value = 44997
$ws.Cells.Item(15,2).Value = "Penelope Fields"
$ws.Cells.Item(15,3).Value = 16
$ws.Cells.Item(15,4).Value = "F"
$ws.Cells.Item(15,5).Value = "Mexe"
$ws.Cells.Item(15,6).Value = 44392
$ws.Cells.Item(15,7).Value = 44393
$ws.Cells.Item(15,11).Value = 44392.66666666667

Write-Output "Done"